# Actualizar grupos de clientes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new client ids to the existing lists in column B
$ws.Range("B10").Value = $ws.Range("B10").Value2 + ".20282"
$ws.Range("B9").Value = $ws.Range("B9").Value2 + ".20205"

# Reflect the selected cell left behind after the edit
$ws.Range("B10").Select()
